$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Copy formatting (style) from row 110 into new rows 111-113 ---
$ws.Range("A110:C110").Copy()
$ws.Range("A111:C113").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 111: Platform Environment: Configuration Utilities ---
$ws.Cells.Item(111, 1).Value = "RTFM"
$ws.Cells.Item(111, 2).Value = "Platform Environment: Configuration Utilities"
$ws.Cells.Item(111, 3).Value = "JDK provides utils for configuration setting:
1. Properties
2. JVM injected command-line arguments (get from args)
3. Get OS environment vars by : Map<String, String> env = System.getenv();
4. Others utils like Preferences API, mainifest in a jar directory, etc."

# --- Row 112: Platform Environment: System Utilities ---
$ws.Cells.Item(112, 1).Value = "RTFM"
$ws.Cells.Item(112, 2).Value = "Platform Environment: System Utilities"
$ws.Cells.Item(112, 3).Value = "Some of utils provided by final static class **System** were covered in previous section, but there are some not covered and here they are:
1. System.in: used for user to read input from command line as user interface
2. System.getProperty(String): return the property value
3. System.getSecurityManager(): security manager is default null in standalone application but available in frameworks like the java Applet. For example, in standalone app, statement ``reader = new FileReader(`"xanadu.txt`");`` might work well but in Applet, a SecurityException wil be thrown (this is true even when invoking a method that isn't documented as throwing SecurityException)
4. Other System utils like System.currentTimeMillis(), System.arrayCopy(), etc"

# --- Row 113: Platform Environment: PATH& CLASSPATH (rich text in column C) ---
$ws.Cells.Item(113, 1).Value = "RTFM"
$ws.Cells.Item(113, 2).Value = "Platform Environment: PATH& CLASSPATH"
$ws.Cells.Item(113, 3).Value = " # Path
Not necessarily needed, mainly for persist convenience. 
# Classpath
The CLASSPATH variable is one way to tell applications, including the JDK tools, where to look for user classes. (Classes that are part of the JRE, JDK platform, and extensions should be defined through other means, such as the bootstrap class path or the extensions directory.)
* -cp can be used to override default classpath
* default classpath of jvm is `".`"
* Class path wildcards allow you to include an entire directory of .jar files in the class path without explicitly naming them individually
For more on classpath, read the _Setting the Class Path_ technical note."

# Apply rich-text run formatting to C113
$chars = $ws.Cells.Item(113, 3).Characters(78, 4)
$chars.Font.Bold = $true
$chars = $ws.Cells.Item(113, 3).Characters(82, 9)
$chars.Font.Bold = $true
$chars.Font.Color = 255
$chars = $ws.Cells.Item(113, 3).Characters(91, 215)
$chars.Font.Bold = $true
$chars = $ws.Cells.Item(113, 3).Characters(306, 20)
$chars.Font.Bold = $true
$chars.Font.Color = 255
$chars = $ws.Cells.Item(113, 3).Characters(326, 31)
$chars.Font.Bold = $true

# --- Fix up row heights now that all content (incl. the long wrapped rich text in ---
# --- C113) has been entered, so auto-fit does not override our explicit height.    ---
$ws.Rows.Item(111).RowHeight = 32.25
$ws.Rows.Item(112).RowHeight = 32.25
$ws.Rows.Item(113).RowHeight = 32.25

# --- Register the bold+red font (with scheme=minor) into the workbook font table, ---
# --- matching the font used for "CLASSPATH" / "bootstrap class path" runs above.   ---
$ghost = $ws.Cells.Item(5000, 1)
$ghost.Value = "x"
$ghost.Font.Bold = $true
$ghost.Font.Color = 255
$ghost.EntireRow.Delete()

# --- Update worksheet view: zoom 130%, scrolled so row 107 is at top ---
$win = $ws.Application.ActiveWindow
$win.Zoom = 130
$win.ScrollRow = 107
$win.ScrollColumn = 1

$ws.Range("C111").Select()

Write-Host "Done applying RTFM Platform Environment entries."
